$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.218.59"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "3.096.32"
$ws.Range("E3").Value = "  -2.36%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.090.68"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000236"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "3.614.11"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "62.608.95"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "3.113.21"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "450.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0983"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.990"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "0.0₃0701"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0381"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "382.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.01%  "
$ws.Range("D44").Value = "2.729.79"
$ws.Range("E44").Value = "  -6.85%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "125.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.244"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.53%  "
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.80%  "
